$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 390 (existing rows 390..462 shift down to 391..463)
$ws.Rows(390).Insert()

# Populate the newly inserted row 390 with the new record
$ws.Range("A390").Value = 10
$ws.Range("B390").Value = "Vega Modelo de Temuco"
$ws.Range("C390").Value = "La Araucanía"
$ws.Range("D390").Value = 45258
$ws.Range("E390").Value = 9
$ws.Range("F390").Value = "Fruta"
$ws.Range("G390").Value = 100103
$ws.Range("H390").Value = "Frutos de hueso (carozo)"
$ws.Range("I390").Value = 100103004
$ws.Range("J390").Value = "Durazno"
$ws.Range("K390").Value = "Florida King"
$ws.Range("L390").Value = "Primera"
$ws.Range("M390").Value = 55
$ws.Range("N390").Value = 25000
$ws.Range("O390").Value = 25000
$ws.Range("P390").Value = 25000
$ws.Range("Q390").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R390").Value = "Provincia de Limarí"
$ws.Range("S390").Value = 1389
$ws.Range("T390").Value = 18
